# Apply the "I0 and IF added" edit: add two new columns (I and J) to the
# sheet, with header labels in row 1 and numeric data in rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they match the rest of the header row (bold, bordered,
# centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-36) -------------------------------------------------
$iValues = @{
    2 = 4; 3 = 1; 4 = 2; 5 = 1; 6 = 1; 7 = 9; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
    29 = 1; 30 = 1; 31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 1; 36 = 1
}
$jValues = @{
    2 = 5; 3 = 6; 4 = 5; 5 = 4; 6 = 5; 7 = 9; 8 = 5; 9 = 4; 10 = 7;
    11 = 4; 12 = 5; 13 = 4; 14 = 6; 15 = 4; 16 = 6; 17 = 7; 18 = 4; 19 = 7;
    20 = 6; 21 = 5; 22 = 7; 23 = 6; 24 = 7; 25 = 6; 26 = 8; 27 = 6; 28 = 7;
    29 = 5; 30 = 6; 31 = 6; 32 = 7; 33 = 4; 34 = 6; 35 = 5; 36 = 2
}

for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]   # column I
    $ws.Cells.Item($row, 10).Value = $jValues[$row]  # column J
}
